# RPA-112: 4-18 varsel yngre enn 15 dgr
#
# The record in row 2 represents an A-Inntekt (income) lookup that failed.
# Clear out the employer/org info that was previously filled in (orgNr,
# firmaNavn, firmaPostAddresse, firmaPostNr) and flag the row with an
# error code in the _5_ErrorCode column (T2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-stale employer lookup results (columns P:S) for the row.
$ws.Range("P2:S2").Clear()

# Record the failure reason in the error-code column.
$ws.Range("T2").Value = "A-Inntekt Failed"

# Widen column G so the (now more relevant) region/info columns are legible,
# and make the selection/scroll reflect where the reviewer was looking.
$ws.Columns("G").ColumnWidth = 21.67

$ws.Range("L2").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1

# Configure the sheet for printing (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
